# Update "想去人数" (want-to-go count, column F) figures on the 展览
# (Exhibition) and 全部类型 (All Types) sheets to match the latest scrape.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 258
$ws1.Range("F8").Value  = 1446
$ws1.Range("F9").Value  = 37829
$ws1.Range("F10").Value = 7752
$ws1.Range("F12").Value = 455
$ws1.Range("F18").Value = 541
$ws1.Range("F20").Value = 67
$ws1.Range("F24").Value = 26
$ws1.Range("F25").Value = 470
$ws1.Range("F27").Value = 466
$ws1.Range("F36").Value = 852

# --- 全部类型 sheet ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 0
$ws4.Range("F4").Value  = 258
$ws4.Range("F9").Value  = 1446
$ws4.Range("F10").Value = 37829
$ws4.Range("F17").Value = 7752
$ws4.Range("F19").Value = 455
$ws4.Range("F26").Value = 541
$ws4.Range("F32").Value = 26
$ws4.Range("F33").Value = 470
$ws4.Range("F35").Value = 466
